# Mise à jour de l'application
# Appends the 2025-08-27 ("45896") training-log entries (rows 205-222)
# to the bottom of the Wellness tracking sheet, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateSerial = 45896

# Row data as it appears in the diff: Name, Volume(C), Intensite(D), Fatigue(E),
# Douleur(F), Localisation douleur(G, optional), Plaisir(H)
$rows = @(
    @{ B = "Emmanuel Valey";   C = 75; D = 4; E = 0; F = 4; G = "Adducteur";  H = 4 },
    @{ B = "Rayane Chayebi";   C = 75; D = 5; E = 7; F = 2; G = "Adducteur "; H = 5 },
    @{ B = "Amir Etien";       C = 75; D = 6; E = 6; F = 0; G = "";           H = 5 },
    @{ B = "Ilyes Boughanmi";  C = 75; D = 5; E = 5; F = 0; G = "";           H = 5 },
    @{ B = "Omar Benyounes";   C = 75; D = 6; E = 6; F = 0; G = "";           H = 5 },
    @{ B = "Yanis Berrached";  C = 75; D = 6; E = 8; F = 0; G = "";           H = 8 },
    @{ B = "Maé Clavel";       C = 75; D = 7; E = 8; F = 7; G = "Ischio";     H = 5 },
    @{ B = "Yoan Zouma";       C = 75; D = 5; E = 7; F = 3; G = "Adducteur "; H = 6 },
    @{ B = "Romain Thunet";    C = 75; D = 7; E = 7; F = 3; G = "Orteil";     H = 2 },
    @{ B = "Jeremie Laurent";  C = 75; D = 7; E = 7; F = 1; G = "Adducteur "; H = 7 },
    @{ B = "Ilan Ihaddadene";  C = 75; D = 4; E = 4; F = 0; G = "";           H = 0 },
    @{ B = "Naim Dhib";        C = 75; D = 5; E = 4; F = 1; G = "Courbature"; H = 5 },
    @{ B = "Amir Kherrab";     C = 75; D = 5; E = 6; F = 0; G = "";           H = 5 },
    @{ B = "Levy Ndoutoume";   C = 75; D = 7; E = 7; F = 5; G = "Genou";      H = 4 },
    @{ B = "Wael Fareh";       C = 75; D = 5; E = 1; F = 0; G = "";           H = 5 },
    @{ B = "Amine Taiar";      C = 75; D = 3; E = 5; F = 6; G = "Semelle";    H = 3 },
    @{ B = "Hedi Nasri";       C = 75; D = 4; E = 5; F = 2; G = "Ischio";     H = 6 },
    @{ B = "Sofiane Belle";    C = 75; D = 4; E = 4; F = 0; G = "";           H = 4 }
)

$r = 205
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yyyy"
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    if ($row.G -ne "") {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"
    $r = $r + 1
}

# Restore the view state recorded in the diff (scroll position + selection).
$ws.Range("K217").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 198
$excel.ActiveWindow.ScrollColumn = 1
